$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for this market/category. It belongs
# chronologically before the existing row 92, so insert a fresh row there,
# which pushes the old rows 92:191 down to 93:192 (dimension grows to R192).
$ws.Rows("92:92").Insert()

$ws.Range("A92").Value = 8
$ws.Range("B92").Value = "Terminal La Palmera de La Serena"
$ws.Range("C92").Value = "Coquimbo"
$ws.Range("D92").Value = 44539
$ws.Range("E92").Value = 4
$ws.Range("F92").Value = 100112012
$ws.Range("G92").Value = "Espinaca"
$ws.Range("H92").Value = "Sin especificar"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 3100
$ws.Range("K92").Value = 400
$ws.Range("L92").Value = 500
$ws.Range("M92").Value = 450
$ws.Range("N92").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O92").Value = "Provincia del Elquí"
$ws.Range("P92").Value = 900
$ws.Range("Q92").Value = 0.5
$ws.Range("R92").Value = "Hortaliza"
